$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New column F: "O QUE FAZER?" (what to do) artefact-delivery notes ---
$ws.Range("F1").Value = "O QUE FAZER? "
$ws.Range("F2").Value = "Cronograma v2.0 , EAP v2.0"
$ws.Range("F5").Value = "Planilha de riscos"
$ws.Range("F7").Value = "Modelo do BD"
$ws.Range("F8").Value = "Diagramas"
$ws.Range("F9").Value = "Plano de testes"
$ws.Range("F10").Value = "Casos de uso"
$ws.Range("F11").Value = "Código fonte v1.0"

# Size the new column like the rest of the table
$ws.Columns("F").ColumnWidth = 40.4

# B9 was retyped by the author (same text, refreshed formatting) as part of
# the "Casos de uso v1.0" sprint-planning pass
$ws.Range("B9").Value = "relatório de testes"

# Leave the selection where the author ended up editing
$ws.Range("B8").Select() | Out-Null
